$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Body")

$ws.Range("G9").Value = "Include 850nm LED (no regulator)"
$ws.Range("G13").Value = "Includes: Objective basket if requested through Cairn research"

# CRISP assembly: wavelength swapped from 850 to 940, supplier Cairn research -> ASI
$ws.Range("B87").Value = "DASI/CRISP-940"
$ws.Range("C87").Value = "ASI"

# ImLock path label update
$ws.Range("A124").Value = "microscope-pathA-imageAutofocus (Imlock)"

# Blackfly camera supplier changed from PointGrey -> FLIR
$ws.Range("C126").Value = "FLIR"

# 900nm Dichroic part number + price swapped for ImLock's actual part; now a long-pass filter
$ws.Range("B151").Value = "DMLP900R"
$ws.Range("E151").Value = 298.67
